# Fantasy.xlsx update — refresh "Final" standings (PF) and "Playoffs"
# (R1/R2/R3) numbers, swap the week-6 team order on the Final tab (and
# introduce the new roster name "Baker Mayzen"), then leave the workbook
# with the "Final" tab active/selected (A17) and "Playoffs" pre-selected
# at C17, matching the author's last on-screen state.

$wb = $excel.ActiveWorkbook

$wsFinal    = $wb.Worksheets.Item("Final")
$wsPlayoffs = $wb.Worksheets.Item("Playoffs")

# ---- Final sheet: updated PF totals ----------------------------------
$wsFinal.Cells.Item(2, 4).Value = 1726.44
$wsFinal.Cells.Item(3, 4).Value = 1721.12

# Row 6 / row 7 swap team names (values, not just positions) and get
# fresh PF totals. Row 7's former name becomes a new roster entry.
$wsFinal.Cells.Item(6, 1).Value = "Red Wave, Red Eyes"
$wsFinal.Cells.Item(6, 4).Value = 1898
$wsFinal.Cells.Item(7, 1).Value = "Baker Mayzen"
$wsFinal.Cells.Item(7, 4).Value = 1887

# ---- Playoffs sheet: updated R1/R2/R3 projections ---------------------
function Set-PlayoffRow($row, $value) {
    $wsPlayoffs.Cells.Item($row, 2).Value = $value
    $wsPlayoffs.Cells.Item($row, 3).Value = $value
    $wsPlayoffs.Cells.Item($row, 4).Value = $value
}

Set-PlayoffRow 2  141.80000000000001
Set-PlayoffRow 3  148.1
Set-PlayoffRow 6  122
Set-PlayoffRow 7  126.1
Set-PlayoffRow 8  155.80000000000001
Set-PlayoffRow 9  128.69999999999999
Set-PlayoffRow 10 146.80000000000001
Set-PlayoffRow 12 129.5
Set-PlayoffRow 13 130.1

# ---- UI state: Playoffs selection moves to C17, then Final becomes ---
# ---- the active tab with A17 selected (matches the saved view) -------
[void]$wsPlayoffs.Activate()
[void]$wsPlayoffs.Range("C17").Select()

[void]$wsFinal.Activate()
[void]$wsFinal.Range("A17").Select()
